$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 08:20"

# Row 17: Austria - numeric updates only
$ws.Cells.Item(17, 2).Value = 11171
$ws.Cells.Item(17, 3).Value = 42
$ws.Cells.Item(17, 5).Value = 9264

# Row 38: Pakistan - numeric updates only
$ws.Cells.Item(38, 2).Value = 2450
$ws.Cells.Item(38, 3).Value = 29
$ws.Cells.Item(38, 4).Value = 126
$ws.Cells.Item(38, 5).Value = 2289
$ws.Cells.Item(38, 6).Value = 10

# Row 77: Kazajistan - numeric updates only
$ws.Cells.Item(77, 2).Value = 448
$ws.Cells.Item(77, 3).Value = 13
$ws.Cells.Item(77, 5).Value = 417

# Rows 85-86: Taiwan and Kuwait swap places (Taiwan gets updated data, Kuwait keeps old data)
$ws.Cells.Item(85, 1).Value = "Taiwan"
$ws.Cells.Item(85, 2).Value = 348
$ws.Cells.Item(85, 3).Value = 9
$ws.Cells.Item(85, 4).Value = 50
$ws.Cells.Item(85, 5).Value = 293
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 5

$ws.Cells.Item(86, 1).Value = "Kuwait"
$ws.Cells.Item(86, 2).Value = 342
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 4).Value = 81
$ws.Cells.Item(86, 5).Value = 261
$ws.Cells.Item(86, 6).Value = 15
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0

# Rows 109-113: Georgia inserted before Venezuela, cascading shift down to Consejo Danes
$ws.Cells.Item(109, 1).Value = "Georgia"
$ws.Cells.Item(109, 2).Value = 148
$ws.Cells.Item(109, 3).Value = 14
$ws.Cells.Item(109, 4).Value = 26
$ws.Cells.Item(109, 5).Value = 122
$ws.Cells.Item(109, 6).Value = 6
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 0

$ws.Cells.Item(110, 1).Value = "Venezuela"
$ws.Cells.Item(110, 2).Value = 146
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = 43
$ws.Cells.Item(110, 5).Value = 98
$ws.Cells.Item(110, 6).Value = 6
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 5

$ws.Cells.Item(111, 1).Value = "Montenegro"
$ws.Cells.Item(111, 2).Value = 144
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 142
$ws.Cells.Item(111, 6).Value = 4
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 2

$ws.Cells.Item(112, 1).Value = "Martinica"
$ws.Cells.Item(112, 2).Value = 138
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 27
$ws.Cells.Item(112, 5).Value = 108
$ws.Cells.Item(112, 6).Value = 19
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 3

$ws.Cells.Item(113, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(113, 2).Value = 134
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 3
$ws.Cells.Item(113, 5).Value = 118
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 13
